$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells that hold numeric-looking text stay text,
# matching the workbook author's convention (inline strings) instead of
# being auto-converted to numbers by Excel's smart entry parsing.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '24.771.89'
$ws.Range("D3").Value = '1.700.53'
$ws.Range("D5").Value = '314.38'
$ws.Range("D7").Value = '0.3980'
$ws.Range("D8").Value = '0.4042'
$ws.Range("D9").Value = '1.001'
$ws.Range("D10").Value = '53.54'
$ws.Range("D11").Value = '1.465'
$ws.Range("D12").Value = '0.08795'
$ws.Range("D13").Value = '26.31'
$ws.Range("D14").Value = '7.518'
$ws.Range("D15").Value = '7.972'
$ws.Range("D16").Value = '0.00001341'
$ws.Range("D17").Value = '1.733.65'
$ws.Range("D18").Value = '95.51'
$ws.Range("D19").Value = '0.07177'
$ws.Range("D21").Value = '7.324'
$ws.Range("D23").Value = '14.38'
$ws.Range("D24").Value = '24.744.73'
$ws.Range("D26").Value = '2.897'
$ws.Range("D27").Value = '23.07'
$ws.Range("D28").Value = '6.117'
$ws.Range("D29").Value = '162.12'
$ws.Range("D30").Value = '144.02'
$ws.Range("D31").Value = '8.296'
$ws.Range("D33").Value = '1.906.99'
$ws.Range("D34").Value = '0.08593'
$ws.Range("D35").Value = '7.321'
$ws.Range("D36").Value = '0.03162'
$ws.Range("D37").Value = '1.025'
$ws.Range("D38").Value = '0.2845'
$ws.Range("D39").Value = '0.09420'
$ws.Range("D41").Value = '10.69'
$ws.Range("D42").Value = '14.14'
$ws.Range("D43").Value = '1.478'
$ws.Range("D44").Value = '17.60'
$ws.Range("D45").Value = '2.699'
$ws.Range("D46").Value = '0.7410'
$ws.Range("D48").Value = '1.375'
$ws.Range("D49").Value = '1.001'
$ws.Range("D50").Value = '0.08378'
$ws.Range("D51").Value = '139.21'

# Restore the default (unformatted) style on the Price column so the
# "@" text-format helper above leaves no visible trace.
$priceRange.Style = "Normal"

# Volume(1h) column (E) values are already non-numeric-looking text
# ("  +0.59%  " style, padded with spaces) so a plain .Value assignment
# keeps them as text without any extra formatting work.
$ws.Range("E2").Value = '  +0.85%  '
$ws.Range("E3").Value = '  +0.70%  '
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("E5").Value = '  +0.10%  '
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("E7").Value = '  +2.43%  '
$ws.Range("E8").Value = '  +0.29%  '
$ws.Range("E9").Value = '  -0.27%  '
$ws.Range("E10").Value = '  +0.94%  '
$ws.Range("E11").Value = '  -1.83%  '
$ws.Range("E12").Value = '  +0.57%  '
$ws.Range("E13").Value = '  +3.54%  '
$ws.Range("E14").Value = '  +0.11%  '
$ws.Range("E15").Value = '  +0.45%  '
$ws.Range("E16").Value = '  -0.95%  '
$ws.Range("E17").Value = '  +2.29%  '
$ws.Range("E18").Value = '  -3.03%  '
$ws.Range("E19").Value = '  +1.16%  '
$ws.Range("E20").Value = '  +4.49%  '
$ws.Range("E21").Value = '  +0.96%  '
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("E23").Value = '  +1.20%  '
$ws.Range("E24").Value = '  +0.69%  '
$ws.Range("E25").Value = '  +1.13%  '
$ws.Range("E26").Value = '  -2.48%  '
$ws.Range("E27").Value = '  +1.66%  '
$ws.Range("E28").Value = '  +17.10%  '
$ws.Range("E29").Value = '  +0.25%  '
$ws.Range("E30").Value = '  +5.25%  '
$ws.Range("E31").Value = '  -5.26%  '
$ws.Range("E32").Value = '  +15.03%  '
$ws.Range("E33").Value = '  +1.47%  '
$ws.Range("E34").Value = '  -2.37%  '
$ws.Range("E35").Value = '  -0.76%  '
$ws.Range("E36").Value = '  +8.66%  '
$ws.Range("E37").Value = '  -0.68%  '
$ws.Range("E38").Value = '  +3.71%  '
$ws.Range("E39").Value = '  +3.38%  '
$ws.Range("E40").Value = '  +5.30%  '
$ws.Range("E41").Value = '  -0.51%  '
$ws.Range("E42").Value = '  -0.63%  '
$ws.Range("E43").Value = '  +1.52%  '
$ws.Range("E44").Value = '  +6.20%  '
$ws.Range("E45").Value = '  +4.37%  '
$ws.Range("E46").Value = '  +3.07%  '
$ws.Range("E47").Value = '  +0.52%  '
$ws.Range("E48").Value = '  +2.83%  '
$ws.Range("E49").Value = '  -0.15%  '
$ws.Range("E50").Value = '  +5.18%  '
$ws.Range("E51").Value = '  +1.02%  '
